# OVW sheet updated with added variation in product listing.
#
# For every locale worksheet in the workbook, cell C5 (the "product listing
# variation" cell) is updated from the shared "product-listing-var3" string
# to a new "product-listing-var7" string, and the sheet's selection is left
# on C5 (mirroring what Excel records when a user clicks/edits that cell).
# The last worksheet processed ends up the active tab, matching the target
# workbook state.

$wb = $excel.ActiveWorkbook
$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    [void]$ws.Activate()
    $ws.Range("C5").Value = "product-listing-var7"
    [void]$ws.Range("C5").Select()
}
